$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4..75 down to 5..76
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new data point
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = 45245
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 100112026
$ws.Cells.Item(4, 7).Value = "Haba"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 13).Value = 10000
$ws.Cells.Item(4, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 400
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
